# Assignment 4: add two new user rows (id 2 and id 4) above the existing
# id=5 row, pushing it from row 3 down to row 5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 3, shifting the current row 3 (id=5) down to row 5.
$ws.Rows("3:4").Insert()

# New row 3: id 2
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "speedywait"
$ws.Range("C3").Value = "ripenmusical"
$ws.Range("D3").Value = "images/pp2.png"

# New row 4: id 4
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "locksignal"
$ws.Range("C4").Value = "opticianblazer"
$ws.Range("D4").Value = "images/pp4.png"
